$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.2721054403052363
$ws.Range("E2").Value = 0.1822636234298907
$ws.Range("F2").Value = 0.5701945524636198

$ws.Range("D3").Value = 0.1617308726398281
$ws.Range("E3").Value = 0.0570192971851197
$ws.Range("F3").Value = 0.201287072166826

$ws.Range("D4").Value = 0.157578289964902
$ws.Range("E4").Value = 0.295338020696204
$ws.Range("F4").Value = 0.654551187720628

$ws.Range("D5").Value = 0.3994996695239591
$ws.Range("E5").Value = 0.2025832219328141
$ws.Range("F5").Value = 0.5690743842060765

$ws.Range("D6").Value = 0.1383861108027632
$ws.Range("E6").Value = 0.06404129843165744
$ws.Range("F6").Value = -0.5344581318440536

$ws.Range("D7").Value = 0.2093966029489203
$ws.Range("E7").Value = 0.1635428229272171
$ws.Range("F7").Value = 0.4356604520789662

$ws.Range("D8").Value = 0.1526217794298124
$ws.Range("E8").Value = 0.0736610121642237
$ws.Range("F8").Value = -0.6451156426589058

$ws.Range("D9").Value = 0.4733114154676658
$ws.Range("E9").Value = 0.2538991046107211
$ws.Range("F9").Value = 0.5008880939628129

$ws.Range("D10").Value = 0.4736425152158926
$ws.Range("E10").Value = 0.3317436295162452
$ws.Range("F10").Value = 0.8667917476135528

$ws.Range("D12").Value = 0.1920887936218409
$ws.Range("E12").Value = 0.1026082095281635
$ws.Range("F12").Value = 0.8434988412851585
